$wb = $excel.ActiveWorkbook

# --- 1. Reorder sheets: move "Eigen_Glasgel" so it sits right after
#        "Eigen_Edelstahl-Horiz" (i.e. before "Eigen_Edelstahl-Stab"). ---
$wsBefore  = $wb.Worksheets.Item("Eigen_Edelstahl-Stab")
$wb.Worksheets.Item("Eigen_Glasgel").Move($wsBefore)

# Re-fetch the handle post-move; the old reference's position is stale.
$wsGlasgel = $wb.Worksheets.Item("Eigen_Glasgel")

# --- 2. Make "Eigen_Glasgel" the active/selected tab. ---
$wsGlasgel.Activate()

# --- 3. Insert a new row for the "Geländer Höhe (m)" (H) variable,
#        right above "Anzahl Ecken". ---
$wsGlasgel.Rows(3).Insert()
$wsGlasgel.Range("A3").Value = "Zahl"
$wsGlasgel.Range("B3").Value = "Geländer Höhe (m)"
$wsGlasgel.Range("C3").Value = "H"
$wsGlasgel.Range("D3").Value = 0

# --- 4. Bump the glass price options (Glas-Typ). ---
$wsGlasgel.Range("D5").Value = "VSG 12.76 Matt:270, VSG 12.76 Klar:260, VSG 10.76 Matt:240, VSG 10.76 Klar:230"

# --- 5. Update the total-price formula to use H and the new constants. ---
$wsGlasgel.Range("E13").Value = "(max(L, 1.0) * H * P_Glas) + (N_Steher * (P_Steher + F_Montage)) + (max(L, 1.0) * P_Handlauf) + ((N_Felder * 4 * P_Klem) + (Ecken * 4 * 20.0)) + (max(L, 1.0) * 75)"

$wsGlasgel.Range("D12").Select()
